# Commit: "Changes made finally. - Commented GitLatch Commit @ 2022-12-28-13-36-21-534"
#
# The underlying edit is the user closing/hiding the GitLatch add-in's task
# pane (word/webextensions/taskpanes.xml: visibility="1" row="1" ->
# visibility="0" row="0"), which in turn makes Word re-stamp the add-in's
# webextension part (word/webextensions/webextension1.xml: id + storeType)
# and drop an unused namespace declaration from word/document.xml when it
# re-serializes the package. All of that lives in add-in/task-pane
# infrastructure parts, not in document body content.
#
# The real Word object model exposes exactly one lever for "task pane
# visibility": Application.TaskPanes(i).Visible. Drive that - it is the
# idiomatic COM-interop equivalent of what the author did in the UI
# (toggling the task pane closed before saving).

$d = $word.ActiveDocument

$paneCount = 0
try {
    $paneCount = $word.TaskPanes.Count
} catch {
    $paneCount = 0
}

if ($paneCount -gt 0) {
    for ($i = 1; $i -le $paneCount; $i++) {
        $pane = $word.TaskPanes.Item($i)
        if ($pane) {
            $pane.Visible = $false
        }
    }
}

$d.Saved = $false
